$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3225
$ws1.Range("F4").Value = 1996
$ws1.Range("F5").Value = 265
$ws1.Range("F6").Value = 102
$ws1.Range("F7").Value = 3079
$ws1.Range("F8").Value = 612
$ws1.Range("F15").Value = 10127
$ws1.Range("F17").Value = 235
$ws1.Range("F20").Value = 8014
$ws1.Range("F21").Value = 12622
$ws1.Range("F24").Value = 20
$ws1.Range("F25").Value = 270
$ws1.Range("F26").Value = 396
$ws1.Range("F28").Value = 6
$ws1.Range("F30").Value = 2819
$ws1.Range("F33").Value = 7944
$ws1.Range("F34").Value = 1495
$ws1.Range("F35").Value = 220
$ws1.Range("F37").Value = 84
$ws1.Range("F38").Value = 4618
$ws1.Range("F39").Value = 1400
$ws1.Range("F40").Value = 73
$ws1.Range("F41").Value = 381

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 8
$ws2.Range("F4").Value = 121

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 145
$ws3.Range("F5").Value = 18

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 3225
$ws4.Range("F7").Value = 1996
$ws4.Range("F8").Value = 8
$ws4.Range("F9").Value = 265
$ws4.Range("F10").Value = 18
$ws4.Range("F11").Value = 3079
$ws4.Range("F12").Value = 121
$ws4.Range("F13").Value = 612
$ws4.Range("F16").Value = 151
$ws4.Range("F19").Value = 10127
$ws4.Range("F20").Value = 235
$ws4.Range("F23").Value = 8014
$ws4.Range("F24").Value = 12622
$ws4.Range("F27").Value = 20
$ws4.Range("F28").Value = 270
$ws4.Range("F32").Value = 6
$ws4.Range("F33").Value = 2819
$ws4.Range("F38").Value = 7944
$ws4.Range("F39").Value = 220
$ws4.Range("F41").Value = 84
$ws4.Range("F42").Value = 4618
